# Update "想去人数" (want-to-go count) figures and one event title across
# the workbook's four sheets, as produced by the latest site regeneration.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 26634
$ws.Range("F4").Value  = 588
$ws.Range("F5").Value  = 254
$ws.Range("F6").Value  = 609
$ws.Range("F7").Value  = 176
$ws.Range("F8").Value  = 547
$ws.Range("F10").Value = 358
$ws.Range("F11").Value = 238
$ws.Range("F12").Value = 189
$ws.Range("F13").Value = 49
$ws.Range("F14").Value = 300
$ws.Range("F15").Value = 61
$ws.Range("F16").Value = 416
$ws.Range("F17").Value = 58
$ws.Range("F18").Value = 1534
$ws.Range("F19").Value = 201
$ws.Range("F20").Value = 41
$ws.Range("F21").Value = 435
$ws.Range("F22").Value = 102

# --- Sheet "演出" (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value  = 4510
$ws.Range("F3").Value  = 231
$ws.Range("C16").Value = "广州·触手猴动漫钢琴音乐演奏会  Marasy Piano Live Asia Tour Prelive "
$ws.Range("F16").Value = 59

# --- Sheet "本地生活" (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5057
$ws.Range("F3").Value = 225

# --- Sheet "全部类型" (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 5057
$ws.Range("F4").Value  = 225
$ws.Range("F5").Value  = 26634
$ws.Range("F6").Value  = 588
$ws.Range("F7").Value  = 4510
$ws.Range("F8").Value  = 254
$ws.Range("F9").Value  = 231
$ws.Range("F10").Value = 609
$ws.Range("F13").Value = 176
$ws.Range("F20").Value = 547
$ws.Range("F23").Value = 358
$ws.Range("F24").Value = 238
$ws.Range("F25").Value = 189
$ws.Range("F26").Value = 49
$ws.Range("F28").Value = 300
$ws.Range("F29").Value = 61
$ws.Range("F32").Value = 416
$ws.Range("F33").Value = 58
$ws.Range("C34").Value = "广州·触手猴动漫钢琴音乐演奏会  Marasy Piano Live Asia Tour Prelive "
$ws.Range("F34").Value = 59
$ws.Range("F35").Value = 1534
$ws.Range("F36").Value = 201
$ws.Range("F38").Value = 41
$ws.Range("F39").Value = 435
$ws.Range("F40").Value = 102
